$wb = $excel.ActiveWorkbook

# --- "About" sheet: update the "last updated" style date in C1 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45392

# --- "MCF" sheet: bump every non-zero/non-one capacity factor to 1 ---
$wsMcf = $wb.Worksheets.Item("MCF")
$wsMcf.Range("B2").Value = 1
$wsMcf.Range("B3").Value = 1
$wsMcf.Range("B4").Value = 1
$wsMcf.Range("B6").Value = 1
$wsMcf.Range("B10").Value = 1
$wsMcf.Range("B11").Value = 1
$wsMcf.Range("B12").Value = 1
$wsMcf.Range("B13").Value = 1
$wsMcf.Range("B14").Value = 1
$wsMcf.Range("B16").Value = 1
$wsMcf.Range("B17").Value = 1
$wsMcf.Range("B18").Value = 1

# Update the selection on the MCF sheet to match the saved view state
$wsMcf.Activate()
$wsMcf.Range("B17").Select()
